$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on columns D and E so numeric-looking
# strings (e.g. "59.29", "0.370") are written as text, matching the
# original inlineStr cell type instead of being auto-coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '36.703.16'
$ws.Range('E2').Value = '  +1.67%  '
$ws.Range('D3').Value = '1.965.40'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').Value = '245.13'
$ws.Range('E5').Value = '  +1.01%  '
$ws.Range('D6').Value = '0.623'
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').Value = '59.29'
$ws.Range('E7').Value = '  +1.77%  '
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('D9').Value = '0.370'
$ws.Range('E9').Value = '  -0.71%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = '0.0865'
$ws.Range('E10').Value = '  +9.42%  '
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').Value = '0.104'
$ws.Range('E11').Value = '  +1.28%  '
$ws.Range('B12').Value = 'Avalanche'
$ws.Range('C12').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D12').Value = '22.42'
$ws.Range('E12').Value = '  +1.01%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').Value = '0.839'
$ws.Range('E13').Value = '  -1.31%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = '13.84'
$ws.Range('E14').Value = '  -0.89%  '
$ws.Range('D15').Value = '2.251.41'
$ws.Range('E15').Value = '  -0.32%  '
$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').Value = '5.29'
$ws.Range('E16').Value = '  -2.22%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '1.961.99'
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '36.584.33'
$ws.Range('E18').Value = '  +1.55%  '
$ws.Range('B19').Value = 'Litecoin'
$ws.Range('C19').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D19').Value = '70.37'
$ws.Range('E19').Value = '  -1.09%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0875'
$ws.Range('E20').Value = '  +2.71%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = '231.21'
$ws.Range('E21').Value = '  -2.10%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '5.09'
$ws.Range('E22').Value = '  -2.17%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').Value = '2.48'
$ws.Range('E24').Value = '  -2.14%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '2.32'
$ws.Range('E25').Value = '  +2.00%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = '9.44'
$ws.Range('E26').Value = '  -3.17%  '
$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').Value = '0.139'
$ws.Range('E27').Value = '  +12.61%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').Value = '162.84'
$ws.Range('E28').Value = '  +1.54%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '19.69'
$ws.Range('E29').Value = '  -0.58%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').Value = '0.119'
$ws.Range('E30').Value = '  -0.47%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = '1.20'
$ws.Range('E31').Value = '  +6.03%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '4.75'
$ws.Range('E32').Value = '  -2.07%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '0.0644'
$ws.Range('E33').Value = '  +4.46%  '
$ws.Range('B34').Value = 'THORChain'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D34').Value = '6.48'
$ws.Range('E34').Value = '  +4.15%  '
$ws.Range('D35').Value = '4.33'
$ws.Range('E35').Value = '  -0.80%  '
$ws.Range('B36').Value = 'BinanceUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  -0.33%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').Value = '2.20'
$ws.Range('E37').Value = '  -3.11%  '
$ws.Range('B38').Value = 'WEMIXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D38').Value = '1.77'
$ws.Range('E38').Value = '  -2.38%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = '3.05'
$ws.Range('E39').Value = '  +0.86%  '
$ws.Range('B40').Value = 'Cronos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D40').Value = '0.1000'
$ws.Range('E40').Value = '  +0.96%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').Value = '1.19'
$ws.Range('E41').Value = '  -2.23%  '
$ws.Range('B42').Value = 'HuobiToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D42').Value = '2.87'
$ws.Range('E42').Value = '  +0.21%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '0.0212'
$ws.Range('E43').Value = '  -0.29%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').Value = '16.47'
$ws.Range('E44').Value = '  +3.64%  '
$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D45').Value = '1.05'
$ws.Range('E45').Value = '  -3.48%  '
$ws.Range('D46').Value = '1.363.46'
$ws.Range('E46').Value = '  +1.76%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = '89.05'
$ws.Range('E47').Value = '  -3.27%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').Value = '7.28'
$ws.Range('E48').Value = '  -3.08%  '
$ws.Range('B49').Value = 'MXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D49').Value = '2.83'
$ws.Range('E49').Value = '  +0.16%  '
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D50').Value = '46.24'
$ws.Range('E50').Value = '  +5.00%  '
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').Value = '1.91'
$ws.Range('E51').Value = '  +9.05%  '

# Restore the default (no explicit style) formatting so the cells
# match the original look (no style index) aside from their content.
$ws.Range("D2:E51").Style = "Normal"
